$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row before row 14: this pushes rows 14-30 down to 15-31,
# matching the target layout (row 13 keeps its number, new content goes
# in row 14, and the rest shift by one).
$ws.Rows("14").Insert()

# --- Row 13: tweak existing remark text ("TBS block" -> "a TBS block") ---
$ws.Range("B13").Value = "* Do not use a formula in a cell that may have its position changed after the merge (for example under a TBS block). Otherwise Excel will raise an error message."

# --- Row 14 (new): extra remark about formula locations ---
$ws.Range("B14").Value = "    This is because the location of formulas are saved a second time in another sub-file for the order of evaluation."

# Rows 15 and 16 already contain the right text after the shift:
#   row15 = "* If a formula uses a reference..." (unchanged)
#   row16 = "* You cannot change picture..."      (unchanged)

# --- Pre-create the number/alignment styles in the same order the
#     original workbook uses them (right-aligned label, bold total,
#     bordered score cell), independent of when the text is written.
$ws.Range("D19").HorizontalAlignment = -4152

$ws.Range("E19").Font.Bold = $true
$ws.Range("E19").NumberFormat = "#,##0.0"

$ws.Range("C21").Copy()
$ws.Range("E21").PasteSpecial(-4122)
$ws.Range("E21").NumberFormat = "#,##0.0"
$ws.Range("E21").HorizontalAlignment = -4152

# --- New "Score" column header on row 20 ---
$ws.Range("B20").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").Value = "Score"

# --- New score field placeholder on row 21 ---
$ws.Range("E21").Value = "[a.score;ope=xlsxNum]"

# --- New "Total:" label + sum formula on row 19 (gap row) ---
$ws.Range("D19").Value = "Total:"
$ws.Range("E19").Formula = "=SUM(E21:E2000)"

# --- Update the active selection to match the saved workbook state ---
[void]$ws.Range("E20").Select()
